# Fixes for TMT data
#
# 1. Design sheet: rename the "Set1"/"Set2" run labels to "A"/"B".
# 2. Fractions sheet: fill in the file -> run mapping for the 8 raw files.
# 3. Selection / active-tab bookkeeping: Design sheet keeps its last
#    selected cell on the last data row, and the Fractions sheet becomes
#    the active tab/selected sheet in the saved workbook.
# 4. Minor column-width tweaks that came along with the edits.

$wb = $excel.ActiveWorkbook

$wsParameters = $wb.Worksheets.Item("Parameters")
$wsDesign     = $wb.Worksheets.Item("Design")
$wsFractions  = $wb.Worksheets.Item("Fractions")

# ---------------------------------------------------------------------
# 1. Design sheet - rename Run labels "Set1" -> "A", "Set2" -> "B"
# ---------------------------------------------------------------------
for ($r = 2; $r -le 21; $r++) {
    $cell = $wsDesign.Cells.Item($r, 1)
    $text = $cell.Text
    if ($text -eq "Set1") {
        $cell.Value = "A"
    } elseif ($text -eq "Set2") {
        $cell.Value = "B"
    }
}

# ---------------------------------------------------------------------
# 2. Fractions sheet - fill in Fraction/Run table
# ---------------------------------------------------------------------

# The sheet previously carried ~300 blank placeholder rows below the
# header; trim them away so the used range shrinks back down to the
# actual data (A1:B9).
$wsFractions.Range("A10:B301").EntireRow.Delete()

$fractionFiles = @("file1.raw", "file2.raw", "file3.raw", "file4.raw", "file5.raw", "file6.raw", "file7.raw", "file8.raw")
$fractionRuns  = @("A", "A", "A", "A", "B", "B", "B", "B")

for ($i = 0; $i -lt $fractionFiles.Length; $i++) {
    $row = $i + 2
    $wsFractions.Cells.Item($row, 1).Value = $fractionFiles[$i]
    $wsFractions.Cells.Item($row, 2).Value = $fractionRuns[$i]
}

# ---------------------------------------------------------------------
# 3. Column width tweaks
# ---------------------------------------------------------------------
$wsParameters.Columns.Item(1).ColumnWidth = 23.862962962962968
$wsParameters.Columns.Item(2).ColumnWidth = 10.729629629629667

$wsDesign.Columns.Item(1).ColumnWidth = 10.729629629629667
$wsDesign.Columns.Item(2).ColumnWidth = 10.729629629629667
$wsDesign.Columns.Item(3).ColumnWidth = 10.729629629629667
$wsDesign.Columns.Item(4).ColumnWidth = 7.103703703703707
$wsDesign.Columns.Item(5).ColumnWidth = 10.729629629629667
$wsDesign.Columns.Item(6).ColumnWidth = 10.729629629629667
$wsDesign.Columns.Item(7).ColumnWidth = 10.729629629629667

$wsFractions.Columns.Item(1).ColumnWidth = 10.729629629629667
$wsFractions.Columns.Item(2).ColumnWidth = 10.729629629629667

# ---------------------------------------------------------------------
# 4. Selection bookkeeping - Design sheet's last active cell moves to A21
# ---------------------------------------------------------------------
$wsDesign.Activate()
$wsDesign.Range("A21").Select()

# ---------------------------------------------------------------------
# 5. Fractions sheet becomes the active/selected tab
# ---------------------------------------------------------------------
$wsFractions.Activate()
$wsFractions.Range("A2").Select()
